$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.397.10'
$ws.Range("E2").Value = '  -2.75%  '
$ws.Range("D3").Value = '3.784.48'
$ws.Range("E3").Value = '  -0.77%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '''593.49'
$ws.Range("E5").Value = '  -1.23%  '
$ws.Range("D6").Value = '''165.80'
$ws.Range("E6").Value = '  -3.07%  '
$ws.Range("D7").Value = '3.781.35'
$ws.Range("E7").Value = '  -0.79%  '
$ws.Range("E8").Value = '  +0.15%  '
$ws.Range("D9").Value = '''0.517'
$ws.Range("E9").Value = '  -1.67%  '
$ws.Range("E12").Value = '  -1.19%  '
$ws.Range("E13").Value = '  -4.01%  '
$ws.Range("D14").Value = '''35.77'
$ws.Range("E14").Value = '  -3.18%  '
$ws.Range("D15").Value = '4.408.52'
$ws.Range("E15").Value = '  -0.69%  '
$ws.Range("D16").Value = '3.775.61'
$ws.Range("E16").Value = '  -0.83%  '
$ws.Range("D17").Value = '67.357.85'
$ws.Range("E17").Value = '  -2.65%  '
$ws.Range("D18").Value = '''17.96'
$ws.Range("E18").Value = '  -1.39%  '
$ws.Range("E19").Value = '  -0.14%  '
$ws.Range("E20").Value = '  -1.99%  '
$ws.Range("E21").Value = '  -7.58%  '
$ws.Range("D22").Value = '''457.78'
$ws.Range("E22").Value = '  -2.96%  '
$ws.Range("D23").Value = '''0.700'
$ws.Range("E23").Value = '  -1.41%  '
$ws.Range("D24").Value = '''0.0000150'
$ws.Range("E24").Value = '  +1.09%  '
$ws.Range("D25").Value = '''83.37'
$ws.Range("E25").Value = '  -1.82%  '
$ws.Range("D26").Value = '''11.84'
$ws.Range("E26").Value = '  -3.36%  '
$ws.Range("D27").Value = '''2.13'
$ws.Range("E27").Value = '  -4.94%  '
$ws.Range("E28").Value = '  +0.07%  '
$ws.Range("D29").Value = '''9.95'
$ws.Range("E29").Value = '  -3.41%  '
$ws.Range("E30").Value = '  -2.15%  '
$ws.Range("D31").Value = '''29.79'
$ws.Range("E31").Value = '  -1.97%  '
$ws.Range("E32").Value = '  -2.76%  '
$ws.Range("D33").Value = '''7.18'
$ws.Range("E33").Value = '  -4.21%  '
$ws.Range("D34").Value = '''9.16'
$ws.Range("E34").Value = '  -2.95%  '
$ws.Range("E35").Value = '  -0.09%  '
$ws.Range("D36").Value = '3.737.55'
$ws.Range("E36").Value = '  -0.73%  '
$ws.Range("D37").Value = '''0.0996'
$ws.Range("E37").Value = '  -2.71%  '
$ws.Range("D38").Value = '''0.138'
$ws.Range("E38").Value = '  -1.38%  '
$ws.Range("D39").Value = '''3.27'
$ws.Range("E39").Value = '  -7.59%  '
$ws.Range("D40").Value = '''0.994'
$ws.Range("E40").Value = '  -1.98%  '
$ws.Range("D41").Value = '''5.72'
$ws.Range("E41").Value = '  -3.02%  '
$ws.Range("D42").Value = '''0.999'
$ws.Range("E42").Value = '  +0.11%  '
$ws.Range("D44").Value = '''43.89'
$ws.Range("E44").Value = '  +0.33%  '
$ws.Range("D45").Value = '''0.298'
$ws.Range("E45").Value = '  -4.15%  '
$ws.Range("D46").Value = '''46.93'
$ws.Range("E46").Value = '  +1.72%  '
$ws.Range("D47").Value = '''8.34'
$ws.Range("E47").Value = '  -3.82%  '
$ws.Range("D48").Value = '''147.52'
$ws.Range("E48").Value = '  +1.66%  '
$ws.Range("D49").Value = '''392.77'
$ws.Range("E49").Value = '  -2.47%  '
$ws.Range("D50").Value = '''1.82'
$ws.Range("E50").Value = '  -8.03%  '
$ws.Range("D51").Value = '2.751.19'
$ws.Range("E51").Value = '  +1.89%  '
